$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the only row whose actual data changes (all other rows just keep
# their original values; their shared-string indices shift automatically
# because of the new strings added to the table below).

# D2: "Ngay bat dau" becomes free text "11/09/2022" instead of a real date.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "11/09/2022"

# E2: "Noi dung" becomes the HTML snippet asking for more members.
$ws.Range("E2").Value = "<p>Cần thêm thành viên tham gia</p>"

# F2, G2, H2 (V1, V2, V3) become real boolean FALSE values.
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

# I2: "Tinh trang" switches from "Da duyet" to "Chua duyet".
$ws.Range("I2").Value = "Chưa duyệt"

# J2: "Link" becomes the text "1" instead of the number 1.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "1"
